# Generate Report for Handback
# Updates the localization-status report:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    (this text is shared by the Overview summary columns and each language
#    sheet's Status column, so every occurrence is updated together)
#  - The per-language "Latest Handback DateTime" is refreshed to the time of
#    this handback run
#  - The stale "handback file is not the latest" Error Detail is cleared now
#    that the handback is in sync
#  - Columns that now hold longer/shorter text are resized to fit

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("K2").Value = "2016-08-17 22:48:19"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("K2").Value = "2016-08-17 22:48:26"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
